# Generate Report for Archive
#
# The two files that were "Ready for handoff" (8fc6fbd0...md and
# 9a31e08b...md) have now moved into translation, so their status rows
# change to "In Translation" on every sheet that tracks status
# (Overview, zh-cn, de-de). The 9a31e08b...md row stays "Ready for
# handoff".

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns are B (zh-cn) and C (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B8").Value = "In Translation"
$overview.Range("C8").Value = "In Translation"
$overview.Range("B9").Value = "In Translation"
$overview.Range("C9").Value = "In Translation"

# --- zh-cn sheet: status column is B ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B8").Value = "In Translation"
$zhcn.Range("B9").Value = "In Translation"

# --- de-de sheet: status column is B ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B8").Value = "In Translation"
$dede.Range("B9").Value = "In Translation"
